$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:AG25").ClearContents()

# Row 1
$ws.Range("A1").Value = 'Name(ID)'
$ws.Range("B1").Value = 'sprite_name'
$ws.Range("C1").Value = 'chat_id'
$ws.Range("D1").Value = 'max_hp'
$ws.Range("E1").Value = 'max_shield'
$ws.Range("F1").Value = 'max_stagger'
$ws.Range("G1").Value = 'ATK'
$ws.Range("H1").Value = 'DEF'
$ws.Range("I1").Value = 'SPEED '
$ws.Range("J1").Value = 'ACC (%)'
$ws.Range("K1").Value = 'EVADE'
$ws.Range("L1").Value = 'AI type'
$ws.Range("M1").Value = 'AI parameters'
$ws.Range("N1").Value = 'null wk/res'
$ws.Range("O1").Value = 'agni wk/res'
$ws.Range("P1").Value = 'cryo wk/res'
$ws.Range("Q1").Value = 'bolt wk/res'
$ws.Range("R1").Value = 'SpellGroup 1'
$ws.Range("S1").Value = 'Spell1 (root)'
$ws.Range("T1").Value = 'Spell1 (elem)'
$ws.Range("U1").Value = 'Spell1(style)'
$ws.Range("V1").Value = 'Spell x (root)'
$ws.Range("W1").Value = 'Spell x (elem)'
$ws.Range("X1").Value = 'Spell x (style)'
$ws.Range("Y1").Value = 'Next Group or END'

# Row 2
$ws.Range("A2").Value = 'Slime'
$ws.Range("B2").Value = 'frog_mario'
$ws.Range("C2").Value = 'enemy_general_1'
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.5
$ws.Range("J2").Value = 1.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 'Attacker1'
$ws.Range("M2").Value = 'none'
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 1.5
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 'GROUP/DEFAULT'
$ws.Range("S2").Value = 'sword'
$ws.Range("T2").Value = 'null'
$ws.Range("U2").Value = 'null'
$ws.Range("V2").Value = 'hammer'
$ws.Range("W2").Value = 'null'
$ws.Range("X2").Value = 'null'
$ws.Range("Y2").Value = 'GROUP/HEALTH_LOW'
$ws.Range("Z2").Value = 'quake'
$ws.Range("AA2").Value = 'null'
$ws.Range("AB2").Value = 'null'
$ws.Range("AC2").Value = 'END'

# Row 3
$ws.Range("A3").Value = 'The Evil Eye'
$ws.Range("B3").Value = 'frog_mario'
$ws.Range("C3").Value = 'enemy_general_1'
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1.25
$ws.Range("H3").Value = 0.2
$ws.Range("I3").Value = 0.75
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 'HealthLow1'
$ws.Range("M3").Value = '75,25'
$ws.Range("N3").Value = -1
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 0.75
$ws.Range("Q3").Value = -2
$ws.Range("R3").Value = 'GROUP/DEFAULT'
$ws.Range("S3").Value = 'sword'
$ws.Range("T3").Value = 'null'
$ws.Range("U3").Value = 'null'
$ws.Range("V3").Value = 'lance'
$ws.Range("W3").Value = 'agni'
$ws.Range("X3").Value = 'null'
$ws.Range("Y3").Value = 'GROUP/HEALTH_LOW'
$ws.Range("Z3").Value = 'quake'
$ws.Range("AA3").Value = 'null'
$ws.Range("AB3").Value = 'null'
$ws.Range("AC3").Value = 'END'

# Row 4
$ws.Range("A4").Value = 'Ladon'
$ws.Range("B4").Value = 'spr_bt_ladon'
$ws.Range("C4").Value = 'enemy_general_1'
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 0.95
$ws.Range("H4").Value = -0.2
$ws.Range("I4").Value = 0.3
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 'Attacker1'
$ws.Range("M4").Value = 'none'
$ws.Range("N4").Value = 1.2
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 'GROUP/DEFAULT'
$ws.Range("S4").Value = 'lance'
$ws.Range("T4").Value = 'null'
$ws.Range("U4").Value = 'null'
$ws.Range("V4").Value = 'lance'
$ws.Range("W4").Value = 'null'
$ws.Range("X4").Value = 'aimed'
$ws.Range("Y4").Value = 'GROUP/HEALTH_LOW'
$ws.Range("Z4").Value = 'quake'
$ws.Range("AA4").Value = 'null'
$ws.Range("AB4").Value = 'null'
$ws.Range("AC4").Value = 'END'

# Row 5
$ws.Range("A5").Value = 'Lilim'
$ws.Range("B5").Value = 'spr_bt_lilim'
$ws.Range("C5").Value = 'enemy_general_1'
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 1.25
$ws.Range("H5").Value = 0.2
$ws.Range("I5").Value = 0.75
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 'HealthLow1'
$ws.Range("M5").Value = '75,25'
$ws.Range("N5").Value = -1
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 0.75
$ws.Range("Q5").Value = -2
$ws.Range("R5").Value = 'GROUP/DEFAULT'
$ws.Range("S5").Value = 'sword'
$ws.Range("T5").Value = 'null'
$ws.Range("U5").Value = 'null'
$ws.Range("V5").Value = 'lance'
$ws.Range("W5").Value = 'agni'
$ws.Range("X5").Value = 'null'
$ws.Range("Y5").Value = 'GROUP/HEALTH_LOW'
$ws.Range("Z5").Value = 'quake'
$ws.Range("AA5").Value = 'null'
$ws.Range("AB5").Value = 'null'
$ws.Range("AC5").Value = 'END'

# Row 6
$ws.Range("A6").Value = 'Changeling'
$ws.Range("B6").Value = 'spr_bt_changeling'
$ws.Range("C6").Value = 'changeling_1'
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 'Attacker1'
$ws.Range("M6").Value = 'none'
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 'GROUP/DEFAULT'
$ws.Range("S6").Value = 'sword'
$ws.Range("T6").Value = 'null'
$ws.Range("U6").Value = 'null'
$ws.Range("V6").Value = 'END'

# Row 7
$ws.Range("A7").Value = 'Wechselbalg'
$ws.Range("B7").Value = 'spr_bt_changeling_healer'
$ws.Range("C7").Value = 'changeling_1'
$ws.Range("D7").Value = 80
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 'Attacker1'
$ws.Range("M7").Value = 'none'
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 'GROUP/DEFAULT'
$ws.Range("S7").Value = 'selfcare'
$ws.Range("T7").Value = 'null'
$ws.Range("U7").Value = 'null'
$ws.Range("V7").Value = 'END'

# Row 8
$ws.Range("A8").Value = 'Tanuki'
$ws.Range("B8").Value = 'spr_bt_tanuki'
$ws.Range("C8").Value = 'tanooki_1'
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.75
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 'Attacker1'
$ws.Range("M8").Value = 'none'
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0.5
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 'GROUP/DEFAULT'
$ws.Range("S8").Value = 'lance'
$ws.Range("T8").Value = 'null'
$ws.Range("U8").Value = 'null'
$ws.Range("V8").Value = 'END'

# Row 9
$ws.Range("A9").Value = 'Tanuki2'
$ws.Range("B9").Value = 'spr_bt_tanuki'
$ws.Range("C9").Value = 'tanooki_1'
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0.85
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 15
$ws.Range("L9").Value = 'Attacker1'
$ws.Range("M9").Value = 'none'
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0.5
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 'GROUP/DEFAULT'
$ws.Range("S9").Value = 'sword'
$ws.Range("T9").Value = 'null'
$ws.Range("U9").Value = 'null'
$ws.Range("V9").Value = 'lance'
$ws.Range("W9").Value = 'null'
$ws.Range("X9").Value = 'null'
$ws.Range("Y9").Value = 'END'

# Row 10
$ws.Range("A10").Value = 'Bunbuku'
$ws.Range("B10").Value = 'spr_bt_tanuki_fire'
$ws.Range("C10").Value = 'tanooki_1'
$ws.Range("D10").Value = 50
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 1.1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0.4
$ws.Range("J10").Value = 0.9
$ws.Range("K10").Value = 15
$ws.Range("L10").Value = 'Attacker1'
$ws.Range("M10").Value = 'none'
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = -1
$ws.Range("P10").Value = 2
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 'GROUP/DEFAULT'
$ws.Range("S10").Value = 'stress'
$ws.Range("T10").Value = 'agni'
$ws.Range("U10").Value = 'null'
$ws.Range("V10").Value = 'sword'
$ws.Range("W10").Value = 'agni'
$ws.Range("X10").Value = 'null'
$ws.Range("Y10").Value = 'END'

# Row 11
$ws.Range("A11").Value = 'Ijiraq'
$ws.Range("B11").Value = 'spr_bt_ijiraq'
$ws.Range("C11").Value = 'ijiraq_1'
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1.2
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.2
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 'Attacker1'
$ws.Range("M11").Value = 'none'
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 2
$ws.Range("P11").Value = -1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 'GROUP/DEFAULT'
$ws.Range("S11").Value = 'sword'
$ws.Range("T11").Value = 'null'
$ws.Range("U11").Value = 'null'
$ws.Range("V11").Value = 'sword'
$ws.Range("W11").Value = 'cryo'
$ws.Range("X11").Value = 'null'
$ws.Range("Y11").Value = 'END'

# Row 12
$ws.Range("A12").Value = 'Ijiraq2'
$ws.Range("B12").Value = 'spr_bt_ijiraq'
$ws.Range("C12").Value = 'ijiraq_1'
$ws.Range("D12").Value = 65
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 1.2
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0.2
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 'Attacker1'
$ws.Range("M12").Value = 'none'
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 2
$ws.Range("P12").Value = -1
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 'GROUP/DEFAULT'
$ws.Range("S12").Value = 'lance'
$ws.Range("T12").Value = 'cryo'
$ws.Range("U12").Value = 'null'
$ws.Range("V12").Value = 'sword'
$ws.Range("W12").Value = 'cryo'
$ws.Range("X12").Value = 'null'
$ws.Range("Y12").Value = 'END'

# Row 13
$ws.Range("A13").Value = 'Ijiraq3'
$ws.Range("B13").Value = 'spr_bt_ijiraq'
$ws.Range("C13").Value = 'ijiraq_1'
$ws.Range("D13").Value = 45
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 0.75
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0.3
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 'Attacker1'
$ws.Range("M13").Value = 'none'
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 1.1
$ws.Range("P13").Value = -1
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 'GROUP/DEFAULT'
$ws.Range("S13").Value = 'lance '
$ws.Range("T13").Value = 'cryo'
$ws.Range("U13").Value = 'null'
$ws.Range("V13").Value = 'sword'
$ws.Range("W13").Value = 'cryo'
$ws.Range("X13").Value = 'null'
$ws.Range("Y13").Value = 'END'

# Row 14
$ws.Range("A14").Value = 'Doppelganger (BLUE)'
$ws.Range("B14").Value = 'spr_bt_doppelganger_b_placeholder'
$ws.Range("C14").Value = 'doppelganger_1'
$ws.Range("D14").Value = 150
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.8
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 'Doppleganger1'
$ws.Range("M14").Value = 'none'
$ws.Range("N14").Value = 0.5
$ws.Range("O14").Value = 2
$ws.Range("P14").Value = -1
$ws.Range("Q14").Value = 0.5
$ws.Range("R14").Value = 'GROUP/DEFAULT'
$ws.Range("S14").Value = 'sword'
$ws.Range("T14").Value = 'cryo'
$ws.Range("U14").Value = 'null'
$ws.Range("V14").Value = 'lance'
$ws.Range("W14").Value = 'cryo'
$ws.Range("X14").Value = 'null'
$ws.Range("Y14").Value = 'GROUP/TOO_LONG'
$ws.Range("Z14").Value = 'hammer'
$ws.Range("AA14").Value = 'cryo'
$ws.Range("AB14").Value = 'null'
$ws.Range("AC14").Value = 'GROUP/SPECIAL'
$ws.Range("AD14").Value = 'magic_circle'
$ws.Range("AE14").Value = 'null'
$ws.Range("AF14").Value = 'null'
$ws.Range("AG14").Value = 'END'

# Row 15
$ws.Range("A15").Value = 'Doppelganger (YELLOW)'
$ws.Range("B15").Value = 'spr_bt_doppelganger_y_placeholder'
$ws.Range("C15").Value = 'doppelganger_1'
$ws.Range("D15").Value = 150
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1.5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.8
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 'Doppleganger1'
$ws.Range("M15").Value = 'none'
$ws.Range("N15").Value = 0.5
$ws.Range("O15").Value = 0.5
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = -1
$ws.Range("R15").Value = 'GROUP/DEFAULT'
$ws.Range("S15").Value = 'sword'
$ws.Range("T15").Value = 'veld'
$ws.Range("U15").Value = 'null'
$ws.Range("V15").Value = 'lance'
$ws.Range("W15").Value = 'veld'
$ws.Range("X15").Value = 'null'
$ws.Range("Y15").Value = 'GROUP/TOO_LONG'
$ws.Range("Z15").Value = 'hammer'
$ws.Range("AA15").Value = 'veld'
$ws.Range("AB15").Value = 'null'
$ws.Range("AC15").Value = 'GROUP/SPECIAL'
$ws.Range("AD15").Value = 'magic_circle'
$ws.Range("AE15").Value = 'null'
$ws.Range("AF15").Value = 'null'
$ws.Range("AG15").Value = 'END'

# Row 16
$ws.Range("A16").Value = 'Doppelganger (RED)'
$ws.Range("B16").Value = 'spr_bt_doppelganger_r_placeholder'
$ws.Range("C16").Value = 'doppelganger_1'
$ws.Range("D16").Value = 150
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 1.5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0.8
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 'Doppleganger1'
$ws.Range("M16").Value = 'none'
$ws.Range("N16").Value = 0.5
$ws.Range("O16").Value = -1
$ws.Range("P16").Value = 0.5
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = 'GROUP/DEFAULT'
$ws.Range("S16").Value = 'sword'
$ws.Range("T16").Value = 'agni'
$ws.Range("U16").Value = 'null'
$ws.Range("V16").Value = 'lance'
$ws.Range("W16").Value = 'agni'
$ws.Range("X16").Value = 'null'
$ws.Range("Y16").Value = 'GROUP/TOO_LONG'
$ws.Range("Z16").Value = 'hammer'
$ws.Range("AA16").Value = 'agni'
$ws.Range("AB16").Value = 'null'
$ws.Range("AC16").Value = 'GROUP/SPECIAL'
$ws.Range("AD16").Value = 'magic_circle'
$ws.Range("AE16").Value = 'null'
$ws.Range("AF16").Value = 'null'
$ws.Range("AG16").Value = 'END'

# Row 17
$ws.Range("A17").Value = 'Doppelganger (???)'
$ws.Range("B17").Value = 'spr_bt_doppelganger_b_placeholder'
$ws.Range("C17").Value = 'doppelganger_1'
$ws.Range("D17").Value = 150
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 1.5
$ws.Range("H17").Value = 0.25
$ws.Range("I17").Value = 0.95
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 5
$ws.Range("L17").Value = 'Doppleganger1'
$ws.Range("M17").Value = 'none'
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 'GROUP/DEFAULT'
$ws.Range("S17").Value = 'sword'
$ws.Range("T17").Value = 'null'
$ws.Range("U17").Value = 'null'
$ws.Range("V17").Value = 'END'

# Row 18
$ws.Range("A18").Value = 'Doppelganger (GRAY)'
$ws.Range("B18").Value = 'spr_bt_doppelganger_g_placeholder'
$ws.Range("C18").Value = 'doppelganger_1'
$ws.Range("D18").Value = 75
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 'Doppleganger1'
$ws.Range("M18").Value = 'none'
$ws.Range("N18").Value = 2
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 'GROUP/DEFAULT'
$ws.Range("S18").Value = 'hammer'
$ws.Range("T18").Value = 'null'
$ws.Range("U18").Value = 'null'
$ws.Range("V18").Value = 'END'

# Row 19
$ws.Range("A19").Value = 'END'

$ws.Range("H9").Select()
